# Auto refresh - 16-02-2026 11:01:35.85
# PaceSmart ML-vs-Excel pacing refresh: recomputed Spend_to_Date / pacing-%
# / Remaining_Budget / Ideal_Daily_Spend figures (and any status flips they
# cause) on Excel_vs_ML, a feature-importance reorder + refresh on
# Feature_Importance, and updated headline KPIs + refresh timestamp on
# Exec_Summary.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Excel_vs_ML - per-campaign recomputed pacing figures
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Excel_vs_ML")

# Row 2
$ws1.Range("H2").Value = 166806.6
$ws1.Range("L2").Value = 166.22
$ws1.Range("M2").Value = -66456.23
# Row 3
$ws1.Range("H3").Value = 89972.7
$ws1.Range("L3").Value = 17.2
$ws1.Range("M3").Value = 462621.93
$ws1.Range("N3").Value = 154207.31
# Row 4
$ws1.Range("H4").Value = 154346.38
$ws1.Range("L4").Value = 52.4
$ws1.Range("M4").Value = 140194
# Row 5
$ws1.Range("H5").Value = 30692.27
$ws1.Range("L5").Value = 6.03
$ws1.Range("M5").Value = 478405.1
# Row 6
$ws1.Range("H6").Value = 146951.89
$ws1.Range("L6").Value = 38.43
$ws1.Range("M6").Value = 235457.36
# Row 7
$ws1.Range("H7").Value = 121020.49
$ws1.Range("L7").Value = 23.02
$ws1.Range("M7").Value = 404675.45
# Row 8
$ws1.Range("H8").Value = 40927.75
$ws1.Range("L8").Value = 70.59
$ws1.Range("M8").Value = 95913.17
$ws1.Range("N8").Value = 2820.98
# Row 10
$ws1.Range("H10").Value = 33037.27
$ws1.Range("L10").Value = 19.07
$ws1.Range("M10").Value = 415820.88
$ws1.Range("N10").Value = 11880.6
# Row 12
$ws1.Range("H12").Value = 86475.62
$ws1.Range("L12").Value = 28.15
$ws1.Range("M12").Value = 473331.11
$ws1.Range("N12").Value = 12792.73
# Row 13
$ws1.Range("H13").Value = 63349.43
$ws1.Range("L13").Value = 97.73
$ws1.Range("M13").Value = 106028.11
$ws1.Range("N13").Value = 2120.56
$ws1.Range("O13").Value = "On Track"
$ws1.Range("V13").Value = "YES"
# Row 14
$ws1.Range("H14").Value = 81515.34
$ws1.Range("L14").Value = 24.34
$ws1.Range("M14").Value = 253420.75
# Row 16
$ws1.Range("H16").Value = 65180.2
$ws1.Range("L16").Value = 40.33
$ws1.Range("M16").Value = 96424.89999999999
# Row 17
$ws1.Range("H17").Value = 139657.66
$ws1.Range("L17").Value = 146.11
$ws1.Range("M17").Value = -44073.85
# Row 20
$ws1.Range("H20").Value = 53250.78
$ws1.Range("L20").Value = 24.59
$ws1.Range("M20").Value = 370429.13
$ws1.Range("N20").Value = 8418.84
# Row 23
$ws1.Range("H23").Value = 150306.27
$ws1.Range("L23").Value = 141.86
$ws1.Range("M23").Value = -44355.2
$ws1.Range("O23").Value = "Overpacing"
# Row 24
$ws1.Range("H24").Value = 118204.01
$ws1.Range("L24").Value = 42.18
$ws1.Range("M24").Value = 162044.58
# Row 25
$ws1.Range("H25").Value = 156018.39
$ws1.Range("L25").Value = 54.58
$ws1.Range("M25").Value = 129844.84
# Row 26
$ws1.Range("H26").Value = 94524.77
$ws1.Range("L26").Value = 42.74
$ws1.Range("M26").Value = 126646.29
# Row 27
$ws1.Range("H27").Value = 85872.27
$ws1.Range("L27").Value = 20.91
$ws1.Range("M27").Value = 491936.61
$ws1.Range("N27").Value = 20497.36
# Row 28
$ws1.Range("H28").Value = 47231.34
$ws1.Range("L28").Value = 24
$ws1.Range("M28").Value = 149546.02
# Row 29
$ws1.Range("H29").Value = 124157.16
$ws1.Range("L29").Value = 72.72
$ws1.Range("M29").Value = 46572.85
# Row 30
$ws1.Range("H30").Value = 83693.11
$ws1.Range("L30").Value = 23.53
$ws1.Range("M30").Value = 271967.15
# Row 31
$ws1.Range("H31").Value = 180010.86
$ws1.Range("L31").Value = 42.98
$ws1.Range("M31").Value = 238846.53
# Row 32
$ws1.Range("H32").Value = 80607.14999999999
$ws1.Range("L32").Value = 23.7
$ws1.Range("M32").Value = 494133.12
$ws1.Range("N32").Value = 10084.35
# Row 33
$ws1.Range("H33").Value = 89176.62
$ws1.Range("L33").Value = 110.97
$ws1.Range("M33").Value = 40180.35
$ws1.Range("N33").Value = 1607.21
$ws1.Range("O33").Value = "Overpacing"
$ws1.Range("S33").Value = "Overdelivered"
$ws1.Range("U33").Value = 40180.35000000001
# Row 34
$ws1.Range("H34").Value = 23560.62
$ws1.Range("L34").Value = 3.98
$ws1.Range("M34").Value = 568389.87
# Row 35
$ws1.Range("H35").Value = 84289.87
$ws1.Range("L35").Value = 21.59
$ws1.Range("M35").Value = 306091.94
# Row 36
$ws1.Range("H36").Value = 188100.04
$ws1.Range("L36").Value = 146.65
$ws1.Range("M36").Value = -59831.25
$ws1.Range("O36").Value = "Overpacing"
# Row 37
$ws1.Range("H37").Value = 3192.36
$ws1.Range("L37").Value = 5.98
$ws1.Range("M37").Value = 530832.14
$ws1.Range("N37").Value = 19660.45
# Row 38
$ws1.Range("H38").Value = 71566.28999999999
$ws1.Range("L38").Value = 18.29
$ws1.Range("M38").Value = 319630.28
# Row 39
$ws1.Range("H39").Value = 81468.39
$ws1.Range("L39").Value = 13.66
$ws1.Range("M39").Value = 514910.54
# Row 40
$ws1.Range("H40").Value = 65098.05
$ws1.Range("L40").Value = 76.59999999999999
$ws1.Range("M40").Value = 54991.4
$ws1.Range("N40").Value = 2894.28
# Row 42
$ws1.Range("H42").Value = 4714.41
$ws1.Range("L42").Value = 9.65
$ws1.Range("M42").Value = 569553.02
$ws1.Range("N42").Value = 13245.42
# Row 43
$ws1.Range("H43").Value = 243342.65
$ws1.Range("L43").Value = 69.08
$ws1.Range("M43").Value = 108914.81
# Row 45
$ws1.Range("H45").Value = 128081.01
$ws1.Range("L45").Value = 43.9
$ws1.Range("M45").Value = 226530
$ws1.Range("N45").Value = 16180.71

# ---------------------------------------------------------------------
# Feature_Importance - reordered + refreshed importances
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("Feature_Importance")

# Row 2
$ws2.Range("A2").Value = "Total_Budget"
$ws2.Range("B2").Value = 0.30403800197451
# Row 3
$ws2.Range("A3").Value = "Pace_Ratio"
$ws2.Range("B3").Value = 0.2978052316094644
# Row 4
$ws2.Range("A4").Value = "Spend_Velocity"
$ws2.Range("B4").Value = 0.1651920240253574
# Row 5
$ws2.Range("A5").Value = "Spend_to_Date"
$ws2.Range("B5").Value = 0.1346997111256371
# Row 6
$ws2.Range("A6").Value = "Days_Elapsed"
$ws2.Range("B6").Value = 0.04185972267982849
# Row 7
$ws2.Range("A7").Value = "Flight_Days"
$ws2.Range("B7").Value = 0.0394229317615561
# Row 8
$ws2.Range("B8").Value = 0.01698237682364667

# ---------------------------------------------------------------------
# Exec_Summary - headline KPIs + refresh timestamp
# ---------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("Exec_Summary")

$ws3.Range("B2").Value = 1
$ws3.Range("B3").Value = 40180.35
$ws3.Range("B5").Value = "2026-02-16 05:31 UTC"

